# Add a new worksheet "Shezwan House 003" (3rd tab) to the Cash Invoice
# workbook, populate it with the new invoice table, and make it the
# active sheet (matching the upstream commit's xlsx diff).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new worksheet after the existing two tabs.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "Shezwan House 003"
$ws3.Tab.Color = 5296274   # RGB(0x92, 0xD0, 0x50) -> FF92D050

# ---------------------------------------------------------------------
# 2. Column widths (matches sheet1 / sheet2 layout conventions).
# ---------------------------------------------------------------------
$ws3.Columns.Item(1).ColumnWidth = 8.88671875
$ws3.Columns.Item(2).ColumnWidth = 48.44140625

# ---------------------------------------------------------------------
# 3. Header row.
# ---------------------------------------------------------------------
$ws3.Range("A1").Value = "SR NO"
$ws3.Range("B1").Value = "ITEM DESCRIPTION"
$ws3.Range("C1").Value = "QTY"
$ws3.Range("D1").Value = "PRICE"
$ws3.Range("E1").Value = "AMOUNT"

# ---------------------------------------------------------------------
# 4. Item rows (2-9).
# ---------------------------------------------------------------------
$items = @(
    @(1, "W Box  2MP Dome Camera", 5, 2000),
    @(2, "BNC Connector", 12, 60),
    @(3, "Power Connector", 6, 40),
    @(4, "Power Supply", 1, 1990),
    @(5, "W Box 2Amps adaptor", 1, 1800),
    @(6, "Enclosure with mounting", 1, 100),
    @(7, "Cabling RJ 59 copper 3 + 1 Cat 6 with casing, cabling, laying ", 87, 90),
    @(8, "INSTALLATION TESTING COMMISSIONING", 1, 4000)
)

$row = 2
foreach ($item in $items) {
    $ws3.Cells.Item($row, 1).Value = $item[0]
    $ws3.Cells.Item($row, 2).Value = $item[1]
    $ws3.Cells.Item($row, 3).Value = $item[2]
    $ws3.Cells.Item($row, 4).Value = $item[3]
    $row = $row + 1
}

# Row 8 ("Cabling ...") wraps onto two lines in the original sheet.
$ws3.Rows.Item(8).RowHeight = 26.4

# Shared AMOUNT formula for the item rows.
$ws3.Range("E2:E9").Formula = "=C2*D2"

# ---------------------------------------------------------------------
# 5. TOTAL row.
# ---------------------------------------------------------------------
$ws3.Range("A10").Value = "TOTAL"
$ws3.Range("A10:D10").Merge()
$ws3.Range("E10").Formula = "=SUM(E2:E9)"

# ---------------------------------------------------------------------
# 6. Formatting - header row: bold 10pt "Calibri  ", boxed, centred, wrap.
# ---------------------------------------------------------------------
$header = $ws3.Range("A1:E1")
$header.Font.Bold = $true
$header.Font.Size = 10
$header.Font.Name = "Calibri  "
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4108
$header.WrapText = $true

# Item rows - regular 10pt "Calibri  ", boxed, centred, wrap.
$body = $ws3.Range("A2:E9")
$body.Font.Size = 10
$body.Font.Name = "Calibri  "
$body.Borders.LineStyle = 1
$body.HorizontalAlignment = -4108
$body.VerticalAlignment = -4108
$body.WrapText = $true

# TOTAL row - bold 10pt "Calibri  ", centred, wrap; box border around the
# merged A10:D10 + E10 cells (distributed across the underlying cells).
$totalRow = $ws3.Range("A10:E10")
$totalRow.Font.Bold = $true
$totalRow.Font.Size = 10
$totalRow.Font.Name = "Calibri  "
$totalRow.HorizontalAlignment = -4108
$totalRow.VerticalAlignment = -4108
$totalRow.WrapText = $true

$ws3.Range("A10:D10").BorderAround(1)
$ws3.Range("E10").Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 7. Selection / activation state.
#    Sheet2 ("Primos Serene 002") keeps its own remembered selection
#    while losing the active tab; Sheet3 becomes the active tab.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Primos Serene 002")
$ws2.Range("C22").Select()

$ws3.Range("H24").Select()
$ws3.Activate()
